$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use row 736 (last fully-populated existing row) as a style template
# for the new rows: column A gets the date number format (style 3),
# columns B/C get style 4, and D.. onward are left unstyled, matching
# the pattern used throughout the sheet.
$templateRow = 736

# Row 737
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A737:K737").PasteSpecial(-4122)
$ws.Cells.Item(737, 1).Value = 44181
$ws.Cells.Item(737, 2).Value = "IC-Z2"
$ws.Cells.Item(737, 3).Value = 1
$ws.Cells.Item(737, 4).Value = 18
$ws.Cells.Item(737, 5).Value = "GOPR0208.MP4_025.png"
$ws.Cells.Item(737, 6).Value = "Ocean Surgeonfish"
$ws.Cells.Item(737, 7).Value = "Acanthurus"
$ws.Cells.Item(737, 8).Value = "bahianus"
$ws.Cells.Item(737, 9).Value = "Acanthurus bahianus"
$ws.Cells.Item(737, 10).Value = 2

# Row 738
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A738:K738").PasteSpecial(-4122)
$ws.Cells.Item(738, 1).Value = 44181
$ws.Cells.Item(738, 2).Value = "IC-Z2"
$ws.Cells.Item(738, 3).Value = 1
$ws.Cells.Item(738, 4).Value = 18
$ws.Cells.Item(738, 5).Value = "GOPR0208.MP4_025.png"
$ws.Cells.Item(738, 6).Value = "Slippery Dick "
$ws.Cells.Item(738, 7).Value = "Halichoeres"
$ws.Cells.Item(738, 8).Value = "bivittatus"
$ws.Cells.Item(738, 9).Value = "Halichoeres bivittatus"
$ws.Cells.Item(738, 10).Value = 1
$ws.Cells.Item(738, 11).Value = "juvenile"

# Row 739
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A739:K739").PasteSpecial(-4122)
$ws.Cells.Item(739, 1).Value = 44181
$ws.Cells.Item(739, 2).Value = "IC-Z2"
$ws.Cells.Item(739, 3).Value = 1
$ws.Cells.Item(739, 4).Value = 18
$ws.Cells.Item(739, 5).Value = "GOPR0208.MP4_025.png"
$ws.Cells.Item(739, 6).Value = "Foureye Butterflyfish"
$ws.Cells.Item(739, 7).Value = "Chaetodon"
$ws.Cells.Item(739, 8).Value = "capistratus"
$ws.Cells.Item(739, 9).Value = "Chaetodon capistratus"
$ws.Cells.Item(739, 10).Value = 1

# Row 740
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A740:K740").PasteSpecial(-4122)
$ws.Cells.Item(740, 1).Value = 44181
$ws.Cells.Item(740, 2).Value = "IC-Z2"
$ws.Cells.Item(740, 3).Value = 1
$ws.Cells.Item(740, 4).Value = 19
$ws.Cells.Item(740, 5).Value = "GOPR0208.MP4_026.png"
$ws.Cells.Item(740, 6).Value = "Ocean Surgeonfish"
$ws.Cells.Item(740, 7).Value = "Acanthurus"
$ws.Cells.Item(740, 8).Value = "bahianus"
$ws.Cells.Item(740, 9).Value = "Acanthurus bahianus"
$ws.Cells.Item(740, 10).Value = 1

# Row 741
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A741:K741").PasteSpecial(-4122)
$ws.Cells.Item(741, 1).Value = 44181
$ws.Cells.Item(741, 2).Value = "IC-Z2"
$ws.Cells.Item(741, 3).Value = 1
$ws.Cells.Item(741, 4).Value = 19
$ws.Cells.Item(741, 5).Value = "GOPR0208.MP4_026.png"
$ws.Cells.Item(741, 6).Value = "Doctorfish"
$ws.Cells.Item(741, 7).Value = "Acanthurus"
$ws.Cells.Item(741, 8).Value = "cirurgus"
$ws.Cells.Item(741, 9).Value = "Acanthurus cirurgus"
$ws.Cells.Item(741, 10).Value = 1

# Row 742
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A742:K742").PasteSpecial(-4122)
$ws.Cells.Item(742, 1).Value = 44181
$ws.Cells.Item(742, 2).Value = "IC-Z2"
$ws.Cells.Item(742, 3).Value = 1
$ws.Cells.Item(742, 4).Value = 19
$ws.Cells.Item(742, 5).Value = "GOPR0208.MP4_026.png"
$ws.Cells.Item(742, 6).Value = "Foureye Butterflyfish"
$ws.Cells.Item(742, 7).Value = "Chaetodon"
$ws.Cells.Item(742, 8).Value = "capistratus"
$ws.Cells.Item(742, 9).Value = "Chaetodon capistratus"
$ws.Cells.Item(742, 10).Value = 2

# Row 743
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A743:K743").PasteSpecial(-4122)
$ws.Cells.Item(743, 1).Value = 44181
$ws.Cells.Item(743, 2).Value = "IC-Z2"
$ws.Cells.Item(743, 3).Value = 1
$ws.Cells.Item(743, 4).Value = 19
$ws.Cells.Item(743, 5).Value = "GOPR0208.MP4_026.png"
$ws.Cells.Item(743, 6).Value = "Slippery Dick "
$ws.Cells.Item(743, 7).Value = "Helichoeres"
$ws.Cells.Item(743, 8).Value = "bivittatus"
$ws.Cells.Item(743, 9).Value = "Helichoeres bivittatus"
$ws.Cells.Item(743, 10).Value = 1
$ws.Cells.Item(743, 11).Value = "juvenile"

# Row 744
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A744:K744").PasteSpecial(-4122)
$ws.Cells.Item(744, 1).Value = 44181
$ws.Cells.Item(744, 2).Value = "IC-Z2"
$ws.Cells.Item(744, 3).Value = 1
$ws.Cells.Item(744, 4).Value = 19
$ws.Cells.Item(744, 5).Value = "GOPR0208.MP4_026.png"
$ws.Cells.Item(744, 6).Value = "Bluehead"
$ws.Cells.Item(744, 7).Value = "Thalassoma"
$ws.Cells.Item(744, 8).Value = "bifasciatum"
$ws.Cells.Item(744, 9).Value = "Thalassoma bifasciatum"
$ws.Cells.Item(744, 10).Value = 1
$ws.Cells.Item(744, 11).Value = "juvenile"

# Row 745
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A745:K745").PasteSpecial(-4122)
$ws.Cells.Item(745, 1).Value = 44181
$ws.Cells.Item(745, 2).Value = "IC-Z1"
$ws.Cells.Item(745, 3).Value = 2
$ws.Cells.Item(745, 4).Value = 1
$ws.Cells.Item(745, 5).Value = "GOPR0204.MP4_010.png"
$ws.Cells.Item(745, 6).Value = "French Grunt"
$ws.Cells.Item(745, 7).Value = "Haemulon"
$ws.Cells.Item(745, 8).Value = "flavolineatum"
$ws.Cells.Item(745, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(745, 10).Value = 9

# Row 746
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A746:K746").PasteSpecial(-4122)
$ws.Cells.Item(746, 1).Value = 44181
$ws.Cells.Item(746, 2).Value = "IC-Z1"
$ws.Cells.Item(746, 3).Value = 2
$ws.Cells.Item(746, 4).Value = 1
$ws.Cells.Item(746, 5).Value = "GOPR0204.MP4_010.png"
$ws.Cells.Item(746, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(746, 7).Value = "Lutjanus "
$ws.Cells.Item(746, 8).Value = "apodus"
$ws.Cells.Item(746, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(746, 10).Value = 4

# Row 747
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A747:K747").PasteSpecial(-4122)
$ws.Cells.Item(747, 1).Value = 44181
$ws.Cells.Item(747, 2).Value = "IC-Z1"
$ws.Cells.Item(747, 3).Value = 2
$ws.Cells.Item(747, 4).Value = 1
$ws.Cells.Item(747, 5).Value = "GOPR0204.MP4_010.png"
$ws.Cells.Item(747, 6).Value = "Yellowtail Damselfish"
$ws.Cells.Item(747, 7).Value = "Microspathodon"
$ws.Cells.Item(747, 8).Value = "chrysurus"
$ws.Cells.Item(747, 9).Value = "Microspathodon chrysurus"
$ws.Cells.Item(747, 10).Value = 1

# Row 748
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A748:K748").PasteSpecial(-4122)
$ws.Cells.Item(748, 1).Value = 44181
$ws.Cells.Item(748, 2).Value = "IC-Z1"
$ws.Cells.Item(748, 3).Value = 2
$ws.Cells.Item(748, 4).Value = 2
$ws.Cells.Item(748, 5).Value = "GOPR0204.MP4_011.png"
$ws.Cells.Item(748, 6).Value = "French Grunt"
$ws.Cells.Item(748, 7).Value = "Haemulon"
$ws.Cells.Item(748, 8).Value = "flavolineatum"
$ws.Cells.Item(748, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(748, 10).Value = 10

# Row 749
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A749:K749").PasteSpecial(-4122)
$ws.Cells.Item(749, 1).Value = 44181
$ws.Cells.Item(749, 2).Value = "IC-Z1"
$ws.Cells.Item(749, 3).Value = 2
$ws.Cells.Item(749, 4).Value = 2
$ws.Cells.Item(749, 5).Value = "GOPR0204.MP4_011.png"
$ws.Cells.Item(749, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(749, 7).Value = "Lutjanus "
$ws.Cells.Item(749, 8).Value = "apodus"
$ws.Cells.Item(749, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(749, 10).Value = 7

# Row 750
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A750:K750").PasteSpecial(-4122)
$ws.Cells.Item(750, 1).Value = 44181
$ws.Cells.Item(750, 2).Value = "IC-Z1"
$ws.Cells.Item(750, 3).Value = 2
$ws.Cells.Item(750, 4).Value = 2
$ws.Cells.Item(750, 5).Value = "GOPR0204.MP4_011.png"
$ws.Cells.Item(750, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(750, 7).Value = "Haemulon"
$ws.Cells.Item(750, 8).Value = "sciurus"
$ws.Cells.Item(750, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(750, 10).Value = 1

# Row 751
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A751:K751").PasteSpecial(-4122)
$ws.Cells.Item(751, 1).Value = 44181
$ws.Cells.Item(751, 2).Value = "IC-Z1"
$ws.Cells.Item(751, 3).Value = 2
$ws.Cells.Item(751, 4).Value = 2
$ws.Cells.Item(751, 5).Value = "GOPR0204.MP4_011.png"
$ws.Cells.Item(751, 6).Value = "Bluehead"
$ws.Cells.Item(751, 7).Value = "Thalassoma"
$ws.Cells.Item(751, 8).Value = "bifasciatum"
$ws.Cells.Item(751, 9).Value = "Thalassoma bifasciatum"
$ws.Cells.Item(751, 10).Value = 1
$ws.Cells.Item(751, 11).Value = "juvenile"

# Row 752
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A752:K752").PasteSpecial(-4122)
$ws.Cells.Item(752, 1).Value = 44181
$ws.Cells.Item(752, 2).Value = "IC-Z1"
$ws.Cells.Item(752, 3).Value = 2
$ws.Cells.Item(752, 4).Value = 3
$ws.Cells.Item(752, 5).Value = "GOPR0204.MP4_013.png"
$ws.Cells.Item(752, 6).Value = "French Grunt"
$ws.Cells.Item(752, 7).Value = "Haemulon"
$ws.Cells.Item(752, 8).Value = "flavolineatum"
$ws.Cells.Item(752, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(752, 10).Value = 5

# Row 753
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A753:K753").PasteSpecial(-4122)
$ws.Cells.Item(753, 1).Value = 44181
$ws.Cells.Item(753, 2).Value = "IC-Z1"
$ws.Cells.Item(753, 3).Value = 2
$ws.Cells.Item(753, 4).Value = 3
$ws.Cells.Item(753, 5).Value = "GOPR0204.MP4_013.png"
$ws.Cells.Item(753, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(753, 7).Value = "Haemulon"
$ws.Cells.Item(753, 8).Value = "sciurus"
$ws.Cells.Item(753, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(753, 10).Value = 1

# Row 754
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A754:K754").PasteSpecial(-4122)
$ws.Cells.Item(754, 1).Value = 44181
$ws.Cells.Item(754, 2).Value = "IC-Z1"
$ws.Cells.Item(754, 3).Value = 2
$ws.Cells.Item(754, 4).Value = 3
$ws.Cells.Item(754, 5).Value = "GOPR0204.MP4_013.png"
$ws.Cells.Item(754, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(754, 7).Value = "Lutjanus "
$ws.Cells.Item(754, 8).Value = "apodus"
$ws.Cells.Item(754, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(754, 10).Value = 4

# Row 755
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A755:K755").PasteSpecial(-4122)
$ws.Cells.Item(755, 1).Value = 44181
$ws.Cells.Item(755, 2).Value = "IC-Z1"
$ws.Cells.Item(755, 3).Value = 2
$ws.Cells.Item(755, 4).Value = 4
$ws.Cells.Item(755, 5).Value = "GOPR0204.MP4_014.png"
$ws.Cells.Item(755, 6).Value = "French Grunt"
$ws.Cells.Item(755, 7).Value = "Haemulon"
$ws.Cells.Item(755, 8).Value = "flavolineatum"
$ws.Cells.Item(755, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(755, 10).Value = 6

# Row 756
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A756:K756").PasteSpecial(-4122)
$ws.Cells.Item(756, 1).Value = 44181
$ws.Cells.Item(756, 2).Value = "IC-Z1"
$ws.Cells.Item(756, 3).Value = 2
$ws.Cells.Item(756, 4).Value = 4
$ws.Cells.Item(756, 5).Value = "GOPR0204.MP4_014.png"
$ws.Cells.Item(756, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(756, 7).Value = "Lutjanus "
$ws.Cells.Item(756, 8).Value = "apodus"
$ws.Cells.Item(756, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(756, 10).Value = 3

# Row 757
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A757:K757").PasteSpecial(-4122)
$ws.Cells.Item(757, 1).Value = 44181
$ws.Cells.Item(757, 2).Value = "IC-Z1"
$ws.Cells.Item(757, 3).Value = 2
$ws.Cells.Item(757, 4).Value = 4
$ws.Cells.Item(757, 5).Value = "GOPR0204.MP4_014.png"
$ws.Cells.Item(757, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(757, 7).Value = "Haemulon"
$ws.Cells.Item(757, 8).Value = "sciurus"
$ws.Cells.Item(757, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(757, 10).Value = 2

# Row 758
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A758:K758").PasteSpecial(-4122)
$ws.Cells.Item(758, 1).Value = 44181
$ws.Cells.Item(758, 2).Value = "IC-Z1"
$ws.Cells.Item(758, 3).Value = 2
$ws.Cells.Item(758, 4).Value = 5
$ws.Cells.Item(758, 5).Value = "GOPR0204.MP4_015.png"
$ws.Cells.Item(758, 6).Value = "French Grunt"
$ws.Cells.Item(758, 7).Value = "Haemulon"
$ws.Cells.Item(758, 8).Value = "flavolineatum"
$ws.Cells.Item(758, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(758, 10).Value = 3

# Row 759
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A759:K759").PasteSpecial(-4122)
$ws.Cells.Item(759, 1).Value = 44181
$ws.Cells.Item(759, 2).Value = "IC-Z1"
$ws.Cells.Item(759, 3).Value = 2
$ws.Cells.Item(759, 4).Value = 5
$ws.Cells.Item(759, 5).Value = "GOPR0204.MP4_015.png"
$ws.Cells.Item(759, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(759, 7).Value = "Lutjanus "
$ws.Cells.Item(759, 8).Value = "apodus"
$ws.Cells.Item(759, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(759, 10).Value = 2

# Row 760
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A760:K760").PasteSpecial(-4122)
$ws.Cells.Item(760, 1).Value = 44181
$ws.Cells.Item(760, 2).Value = "IC-Z1"
$ws.Cells.Item(760, 3).Value = 2
$ws.Cells.Item(760, 4).Value = 5
$ws.Cells.Item(760, 5).Value = "GOPR0204.MP4_015.png"
$ws.Cells.Item(760, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(760, 7).Value = "Haemulon"
$ws.Cells.Item(760, 8).Value = "sciurus"
$ws.Cells.Item(760, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(760, 10).Value = 1

# Row 761
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A761:K761").PasteSpecial(-4122)
$ws.Cells.Item(761, 1).Value = 44181
$ws.Cells.Item(761, 2).Value = "IC-Z1"
$ws.Cells.Item(761, 3).Value = 2
$ws.Cells.Item(761, 4).Value = 6
$ws.Cells.Item(761, 5).Value = "GOPR0204.MP4_016.png"
$ws.Cells.Item(761, 6).Value = "French Grunt"
$ws.Cells.Item(761, 7).Value = "Haemulon"
$ws.Cells.Item(761, 8).Value = "flavolineatum"
$ws.Cells.Item(761, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(761, 10).Value = 6

# Row 762
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A762:K762").PasteSpecial(-4122)
$ws.Cells.Item(762, 1).Value = 44181
$ws.Cells.Item(762, 2).Value = "IC-Z1"
$ws.Cells.Item(762, 3).Value = 2
$ws.Cells.Item(762, 4).Value = 6
$ws.Cells.Item(762, 5).Value = "GOPR0204.MP4_016.png"
$ws.Cells.Item(762, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(762, 7).Value = "Lutjanus "
$ws.Cells.Item(762, 8).Value = "apodus"
$ws.Cells.Item(762, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(762, 10).Value = 2

# Row 763
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A763:K763").PasteSpecial(-4122)
$ws.Cells.Item(763, 1).Value = 44181
$ws.Cells.Item(763, 2).Value = "IC-Z1"
$ws.Cells.Item(763, 3).Value = 2
$ws.Cells.Item(763, 4).Value = 6
$ws.Cells.Item(763, 5).Value = "GOPR0204.MP4_016.png"
$ws.Cells.Item(763, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(763, 7).Value = "Haemulon"
$ws.Cells.Item(763, 8).Value = "sciurus"
$ws.Cells.Item(763, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(763, 10).Value = 4

# Row 764
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A764:K764").PasteSpecial(-4122)
$ws.Cells.Item(764, 1).Value = 44181
$ws.Cells.Item(764, 2).Value = "IC-Z1"
$ws.Cells.Item(764, 3).Value = 2
$ws.Cells.Item(764, 4).Value = 6
$ws.Cells.Item(764, 5).Value = "GOPR0204.MP4_016.png"
$ws.Cells.Item(764, 6).Value = "Bicolor Damselfish"
$ws.Cells.Item(764, 7).Value = "Stegastes"
$ws.Cells.Item(764, 8).Value = "partitus"
$ws.Cells.Item(764, 9).Value = "Stegastes partitus"
$ws.Cells.Item(764, 10).Value = 2

# Row 765
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A765:K765").PasteSpecial(-4122)
$ws.Cells.Item(765, 1).Value = 44181
$ws.Cells.Item(765, 2).Value = "IC-Z1"
$ws.Cells.Item(765, 3).Value = 2
$ws.Cells.Item(765, 4).Value = 7
$ws.Cells.Item(765, 5).Value = "GOPR0204.MP4_017.png"
$ws.Cells.Item(765, 6).Value = "French Grunt"
$ws.Cells.Item(765, 7).Value = "Haemulon"
$ws.Cells.Item(765, 8).Value = "flavolineatum"
$ws.Cells.Item(765, 9).Value = "Haemulon flavolineatum"
$ws.Cells.Item(765, 10).Value = 3

# Row 766
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A766:K766").PasteSpecial(-4122)
$ws.Cells.Item(766, 1).Value = 44181
$ws.Cells.Item(766, 2).Value = "IC-Z1"
$ws.Cells.Item(766, 3).Value = 2
$ws.Cells.Item(766, 4).Value = 7
$ws.Cells.Item(766, 5).Value = "GOPR0204.MP4_017.png"
$ws.Cells.Item(766, 6).Value = "Schoolmaster Snapper"
$ws.Cells.Item(766, 7).Value = "Lutjanus "
$ws.Cells.Item(766, 8).Value = "apodus"
$ws.Cells.Item(766, 9).Value = "Lutjanus apodus"
$ws.Cells.Item(766, 10).Value = 2

# Row 767
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A767:K767").PasteSpecial(-4122)
$ws.Cells.Item(767, 1).Value = 44181
$ws.Cells.Item(767, 2).Value = "IC-Z1"
$ws.Cells.Item(767, 3).Value = 2
$ws.Cells.Item(767, 4).Value = 7
$ws.Cells.Item(767, 5).Value = "GOPR0204.MP4_017.png"
$ws.Cells.Item(767, 6).Value = "Princess Parrotfish"
$ws.Cells.Item(767, 7).Value = "Scarus"
$ws.Cells.Item(767, 8).Value = "taeniopterus"
$ws.Cells.Item(767, 9).Value = "Scarus taeniopterus"
$ws.Cells.Item(767, 10).Value = 1
$ws.Cells.Item(767, 11).Value = "juvenile"

# Row 768
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A768:K768").PasteSpecial(-4122)
$ws.Cells.Item(768, 1).Value = 44181
$ws.Cells.Item(768, 2).Value = "IC-Z1"
$ws.Cells.Item(768, 3).Value = 2
$ws.Cells.Item(768, 4).Value = 7
$ws.Cells.Item(768, 5).Value = "GOPR0204.MP4_017.png"
$ws.Cells.Item(768, 6).Value = "Bicolor Damselfish"
$ws.Cells.Item(768, 7).Value = "Stegates"
$ws.Cells.Item(768, 8).Value = "partitus"
$ws.Cells.Item(768, 9).Value = "Stegastes partitus"
$ws.Cells.Item(768, 10).Value = 1

# Row 769
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A769:K769").PasteSpecial(-4122)
$ws.Cells.Item(769, 1).Value = 44181
$ws.Cells.Item(769, 2).Value = "IC-Z1"
$ws.Cells.Item(769, 3).Value = 2
$ws.Cells.Item(769, 4).Value = 8
$ws.Cells.Item(769, 5).Value = "GOPR0204.MP4_018.png"
$ws.Cells.Item(769, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(769, 7).Value = "Haemulon"
$ws.Cells.Item(769, 8).Value = "sciurus"
$ws.Cells.Item(769, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(769, 10).Value = 2

# Row 770
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A770:K770").PasteSpecial(-4122)
$ws.Cells.Item(770, 1).Value = 44181
$ws.Cells.Item(770, 2).Value = "IC-Z1"
$ws.Cells.Item(770, 3).Value = 2
$ws.Cells.Item(770, 4).Value = 8
$ws.Cells.Item(770, 5).Value = "GOPR0204.MP4_018.png"
$ws.Cells.Item(770, 6).Value = "Princess Parrotfish"
$ws.Cells.Item(770, 7).Value = "Scarus "
$ws.Cells.Item(770, 8).Value = "taeniopterus"
$ws.Cells.Item(770, 9).Value = "Scarus taeniopterus"
$ws.Cells.Item(770, 10).Value = 1
$ws.Cells.Item(770, 11).Value = "juvenile"

# Row 771
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A771:K771").PasteSpecial(-4122)
$ws.Cells.Item(771, 1).Value = 44181
$ws.Cells.Item(771, 2).Value = "IC-Z1"
$ws.Cells.Item(771, 3).Value = 2
$ws.Cells.Item(771, 4).Value = 9
$ws.Cells.Item(771, 5).Value = "GOPR0204.MP4_019.png"
$ws.Cells.Item(771, 6).Value = "Bluestriped Grunt"
$ws.Cells.Item(771, 7).Value = "Haemulon"
$ws.Cells.Item(771, 8).Value = "sciurus"
$ws.Cells.Item(771, 9).Value = "Haemulon sciurus"
$ws.Cells.Item(771, 10).Value = 1

# Row 772
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A772:K772").PasteSpecial(-4122)
$ws.Cells.Item(772, 1).Value = 44181
$ws.Cells.Item(772, 2).Value = "IC-Z1"
$ws.Cells.Item(772, 3).Value = 2
$ws.Cells.Item(772, 4).Value = 9
$ws.Cells.Item(772, 5).Value = "GOPR0204.MP4_019.png"
$ws.Cells.Item(772, 6).Value = "Princess Parrotfish"
$ws.Cells.Item(772, 7).Value = "Scarus "
$ws.Cells.Item(772, 8).Value = "taeniopterus"
$ws.Cells.Item(772, 9).Value = "Scarus taeniopterus"
$ws.Cells.Item(772, 10).Value = 2
$ws.Cells.Item(772, 11).Value = "juvenile"

# Row 773
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A773:K773").PasteSpecial(-4122)
$ws.Cells.Item(773, 1).Value = 44181
$ws.Cells.Item(773, 2).Value = "IC-Z1"
$ws.Cells.Item(773, 3).Value = 2
$ws.Cells.Item(773, 4).Value = 9
$ws.Cells.Item(773, 5).Value = "GOPR0204.MP4_019.png"
$ws.Cells.Item(773, 6).Value = "Yellowhead Wrasse"
$ws.Cells.Item(773, 7).Value = "Halichoeres"
$ws.Cells.Item(773, 8).Value = "garnoti"
$ws.Cells.Item(773, 9).Value = "Halichoeres garnoti"
$ws.Cells.Item(773, 10).Value = 1
$ws.Cells.Item(773, 11).Value = "initial"

# Row 774
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A774:K774").PasteSpecial(-4122)
$ws.Cells.Item(774, 1).Value = 44181
$ws.Cells.Item(774, 2).Value = "IC-Z1"
$ws.Cells.Item(774, 3).Value = 2
$ws.Cells.Item(774, 4).Value = 10
$ws.Cells.Item(774, 5).Value = "GOPR0204.MP4_020.png"
$ws.Cells.Item(774, 6).Value = "Bicolor Damselfish"
$ws.Cells.Item(774, 7).Value = "Stegastes"
$ws.Cells.Item(774, 8).Value = "partitus"
$ws.Cells.Item(774, 9).Value = "Stegastes partitus"
$ws.Cells.Item(774, 10).Value = 1

# Row 775
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A775:K775").PasteSpecial(-4122)
$ws.Cells.Item(775, 1).Value = 44181
$ws.Cells.Item(775, 2).Value = "IC-Z1"
$ws.Cells.Item(775, 3).Value = 2
$ws.Cells.Item(775, 4).Value = 10
$ws.Cells.Item(775, 5).Value = "GOPR0204.MP4_020.png"
$ws.Cells.Item(775, 6).Value = "Foureye Butterflyfish"
$ws.Cells.Item(775, 7).Value = "Chaetodon"
$ws.Cells.Item(775, 8).Value = "capistratus"
$ws.Cells.Item(775, 9).Value = "Chaetodon capistratus"
$ws.Cells.Item(775, 10).Value = 1

# Row 776
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A776:K776").PasteSpecial(-4122)
$ws.Cells.Item(776, 1).Value = 44181
$ws.Cells.Item(776, 2).Value = "IC-Z1"
$ws.Cells.Item(776, 3).Value = 2
$ws.Cells.Item(776, 4).Value = 11
$ws.Cells.Item(776, 5).Value = "GOPR0204.MP4_021.png"
$ws.Cells.Item(776, 6).Value = "Bicolor Damselfish"
$ws.Cells.Item(776, 7).Value = "Stegastes"
$ws.Cells.Item(776, 8).Value = "partitus"
$ws.Cells.Item(776, 9).Value = "Stegastes partitus"
$ws.Cells.Item(776, 10).Value = 1

# Row 777
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A777:K777").PasteSpecial(-4122)
$ws.Cells.Item(777, 1).Value = 44181
$ws.Cells.Item(777, 2).Value = "IC-Z1"
$ws.Cells.Item(777, 3).Value = 2
$ws.Cells.Item(777, 4).Value = 12
$ws.Cells.Item(777, 5).Value = "GOPR0204.MP4_022.png"
$ws.Cells.Item(777, 6).Value = "Bicolor Damselfish"
$ws.Cells.Item(777, 7).Value = "Stegastes"
$ws.Cells.Item(777, 8).Value = "partitus"
$ws.Cells.Item(777, 9).Value = "Stegastes partitus"
$ws.Cells.Item(777, 10).Value = 1

# Row 778
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A778:K778").PasteSpecial(-4122)
$ws.Cells.Item(778, 1).Value = 44181
$ws.Cells.Item(778, 2).Value = "IC-Z1"
$ws.Cells.Item(778, 3).Value = 2
$ws.Cells.Item(778, 4).Value = 12
$ws.Cells.Item(778, 5).Value = "GOPR0204.MP4_022.png"
$ws.Cells.Item(778, 6).Value = "Bluehead"
$ws.Cells.Item(778, 7).Value = "Thalassoma"
$ws.Cells.Item(778, 8).Value = "bifasciatum"
$ws.Cells.Item(778, 9).Value = "Thalassoma bifasciatum"
$ws.Cells.Item(778, 10).Value = 1
$ws.Cells.Item(778, 11).Value = "juvenile"

# Row 779
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A779:K779").PasteSpecial(-4122)
$ws.Cells.Item(779, 1).Value = 44181
$ws.Cells.Item(779, 2).Value = "IC-Z1"
$ws.Cells.Item(779, 3).Value = 2
$ws.Cells.Item(779, 4).Value = 13
$ws.Cells.Item(779, 5).Value = "GOPR0204.MP4_023.png"
$ws.Cells.Item(779, 6).Value = "Bluehead"
$ws.Cells.Item(779, 7).Value = "Thalassoma"
$ws.Cells.Item(779, 8).Value = "bifasciatum"
$ws.Cells.Item(779, 9).Value = "Thalassoma bifasciatum"
$ws.Cells.Item(779, 10).Value = 3
$ws.Cells.Item(779, 11).Value = "juvenile"

# Row 780
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A780:K780").PasteSpecial(-4122)
$ws.Cells.Item(780, 1).Value = 44181
$ws.Cells.Item(780, 2).Value = "IC-Z1"
$ws.Cells.Item(780, 3).Value = 2
$ws.Cells.Item(780, 4).Value = 14
$ws.Cells.Item(780, 5).Value = "GOPR0204.MP4_024.png"
$ws.Cells.Item(780, 6).Value = "Bluehead"
$ws.Cells.Item(780, 7).Value = "Thalassoma"
$ws.Cells.Item(780, 8).Value = "bifasciatum"
$ws.Cells.Item(780, 9).Value = "Thalassoma bifasciatum"
$ws.Cells.Item(780, 10).Value = 1
$ws.Cells.Item(780, 11).Value = "juvenile"

# Row 781
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A781:K781").PasteSpecial(-4122)
$ws.Cells.Item(781, 1).Value = 44181
$ws.Cells.Item(781, 2).Value = "IC-Z1"
$ws.Cells.Item(781, 3).Value = 2
$ws.Cells.Item(781, 4).Value = 15
$ws.Cells.Item(781, 5).Value = "GOPR0204.MP4_025.png"
$ws.Cells.Item(781, 6).Value = "Bicolor Damselfish"
$ws.Cells.Item(781, 7).Value = "Stegastes"
$ws.Cells.Item(781, 8).Value = "partitus"
$ws.Cells.Item(781, 9).Value = "Stegastes partitus"
$ws.Cells.Item(781, 10).Value = 3

# Row 782
$ws.Range("A" + $templateRow + ":K" + $templateRow).Copy()
$ws.Range("A782:K782").PasteSpecial(-4122)
$ws.Cells.Item(782, 1).Value = 44181
$ws.Cells.Item(782, 2).Value = "IC-Z1"
$ws.Cells.Item(782, 3).Value = 2
$ws.Cells.Item(782, 4).Value = 16
$ws.Cells.Item(782, 5).Value = "GOPR0204.MP4_026.png"

$excel.CutCopyMode = 0

# Update the frozen-pane view + selection to match the target state
$ws.Range("F782").Select()
